# Search Product check Table (cột tên)
# Replace the mismatched "ParentCategory" / "Category Name" values in the
# Category and Product sheets with the correct category names, and update
# the active selections left behind on each sheet.

$wb = $excel.ActiveWorkbook

# --- Category sheet --------------------------------------------------
# ParentCategory column (B) incorrectly held product-ish names; normalise
# every row to the single real parent category.
$wsCategory = $wb.Worksheets.Item("Category")
$wsCategory.Activate()
$wsCategory.Range("B2").Value = "Category 1"
$wsCategory.Range("B3").Value = "Category 1"
$wsCategory.Range("B4").Value = "Category 1"
$wsCategory.Range("B2").Select()

# --- Product sheet -----------------------------------------------------
# Category Name column (B) pointed at brand/old placeholder text instead
# of real category names.
$wsProduct = $wb.Worksheets.Item("Product")
$wsProduct.Activate()
$wsProduct.Range("B2").Value = "Mỹ Phẩm"
$wsProduct.Range("B3").Value = "Cake"
$wsProduct.Range("I2").Select()

# --- Brand sheet ---------------------------------------------------------
# Just a leftover selection change, no data edits.
$wsBrand = $wb.Worksheets.Item("Brand")
$wsBrand.Activate()
$wsBrand.Range("C1").Select()

# Restore focus to the Product sheet, which remains the active tab.
$wsProduct.Activate()
